$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Per-row updates derived from the source diff (Coin=B, Link=C, Price=D, Volume(1h)=E)
$updates = @(
    [PSCustomObject]@{ Row = 2; D = "59.931.44"; E = "  +0.33%  " },
    [PSCustomObject]@{ Row = 3; D = "2.538.83"; E = "  +0.44%  " },
    [PSCustomObject]@{ Row = 4; E = "  -0.01%  " },
    [PSCustomObject]@{ Row = 5; D = "544.49"; E = "  +0.29%  " },
    [PSCustomObject]@{ Row = 6; D = "145.37"; E = "  -1.69%  " },
    [PSCustomObject]@{ Row = 7; E = "  -0.42%  " },
    [PSCustomObject]@{ Row = 8; D = "0.574"; E = "  -1.32%  " },
    [PSCustomObject]@{ Row = 9; D = "2.570.59"; E = "  +1.80%  " },
    [PSCustomObject]@{ Row = 10; D = "0.101"; E = "  +0.34%  " },
    [PSCustomObject]@{ Row = 11; E = "  +1.21%  " },
    [PSCustomObject]@{ Row = 12; D = "5.58"; E = "  +2.61%  " },
    [PSCustomObject]@{ Row = 13; E = "  +1.04%  " },
    [PSCustomObject]@{ Row = 14; D = "2.988.84"; E = "  +0.37%  " },
    [PSCustomObject]@{ Row = 15; D = "23.79"; E = "  -3.55%  " },
    [PSCustomObject]@{ Row = 16; D = "59.872.42"; E = "  +0.07%  " },
    [PSCustomObject]@{ Row = 17; D = "0.0000143"; E = "  +2.11%  " },
    [PSCustomObject]@{ Row = 18; D = "2.554.59"; E = "  +1.79%  " },
    [PSCustomObject]@{ Row = 19; D = "11.30"; E = "  -3.19%  " },
    [PSCustomObject]@{ Row = 20; D = "4.33"; E = "  -1.40%  " },
    [PSCustomObject]@{ Row = 21; D = "328.51"; E = "  +0.29%  " },
    [PSCustomObject]@{ Row = 22; D = "0.998"; E = "  -0.37%  " },
    [PSCustomObject]@{ Row = 23; D = "5.96"; E = "  +2.33%  " },
    [PSCustomObject]@{ Row = 24; D = "62.29"; E = "  +1.26%  " },
    [PSCustomObject]@{ Row = 25; E = "  -1.69%  " },
    [PSCustomObject]@{ Row = 26; D = "0.167"; E = "  +2.41%  " },
    [PSCustomObject]@{ Row = 27; E = "  -1.97%  " },
    [PSCustomObject]@{ Row = 28; D = "8.03"; E = "  +1.14%  " },
    [PSCustomObject]@{ Row = 29; B = "Aptos"; C = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"; D = "7.04"; E = "  -2.77%  " },
    [PSCustomObject]@{ Row = 30; B = "PEPE"; C = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"; D = "0.0₃0799"; E = "  +0.37%  " },
    [PSCustomObject]@{ Row = 31; E = "  -0.98%  " },
    [PSCustomObject]@{ Row = 32; D = "1.22"; E = "  -5.66%  " },
    [PSCustomObject]@{ Row = 33; D = "162.41"; E = "  +1.80%  " },
    [PSCustomObject]@{ Row = 34; B = "USDe"; C = "https://coinranking.com/coin/exbfr2U-0+usde-usde"; D = "0.997"; E = "  -0.16%  " },
    [PSCustomObject]@{ Row = 35; B = "ImmutableX"; C = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"; D = "1.47"; E = "  +2.15%  " },
    [PSCustomObject]@{ Row = 36; D = "18.84"; E = "  +0.44%  " },
    [PSCustomObject]@{ Row = 37; D = "4.48"; E = "  -0.78%  " },
    [PSCustomObject]@{ Row = 38; E = "  -3.50%  " },
    [PSCustomObject]@{ Row = 39; D = "5.69"; E = "  -6.35%  " },
    [PSCustomObject]@{ Row = 40; D = "37.21"; E = "  +1.34%  " },
    [PSCustomObject]@{ Row = 41; D = "302.66"; E = "  -4.66%  " },
    [PSCustomObject]@{ Row = 42; D = "0.840"; E = "  +0.23%  " },
    [PSCustomObject]@{ Row = 43; D = "3.74"; E = "  -1.73%  " },
    [PSCustomObject]@{ Row = 44; B = "Mantle"; C = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"; D = "0.610"; E = "  +1.08%  " },
    [PSCustomObject]@{ Row = 45; B = "FirstDigitalUSD"; C = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"; D = "0.993"; E = "  -0.39%  " },
    [PSCustomObject]@{ Row = 46; D = "10.83"; E = "  +0.28%  " },
    [PSCustomObject]@{ Row = 47; D = "19.04"; E = "  +1.74%  " },
    [PSCustomObject]@{ Row = 48; E = "  -0.77%  " },
    [PSCustomObject]@{ Row = 49; D = "124.90"; E = "  -1.70%  " },
    [PSCustomObject]@{ Row = 50; D = "0.0521"; E = "  -2.17%  " },
    [PSCustomObject]@{ Row = 51; E = "  -1.42%  " },
)

foreach ($u in $updates) {
    $r = $u.Row
    if ($u.B) {
        $ws.Cells.Item($r, 2).Value = $u.B
    }
    if ($u.C) {
        $ws.Cells.Item($r, 3).Value = $u.C
    }
    if ($u.D) {
        # Force text storage so numeric-looking strings (e.g. "5.58") are not
        # reinterpreted as numbers, then restore the default (unstyled) look.
        $cellD = $ws.Cells.Item($r, 4)
        $cellD.NumberFormat = "@"
        $cellD.Value = $u.D
        $cellD.Style = "Normal"
    }
    if ($u.E) {
        $ws.Cells.Item($r, 5).Value = $u.E
    }
}

Write-Host ("Applied " + $updates.Count + " row updates")
